# Auto-generated Excel COM-interop script.
# Applies a scheduled market-data refresh to the Leve profit-tracking sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across all class tabs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6898.6
$ws.Range("I33").Value = 98.42856999999999
$ws.Range("J33").Value = 12848.75
$ws.Range("K33").Value = 98.42856999999999
$ws.Range("L33").Value = 12848.75
$ws.Range("M33").Value = 130.57143
$ws.Range("N33").Value = -13306.75
$ws.Range("H53").Value = 1443.3846
$ws.Range("I53").Value = 1685
$ws.Range("J53").Value = 1161.5
$ws.Range("K53").Value = 1685
$ws.Range("L53").Value = 1161.5
$ws.Range("M53").Value = -1048
$ws.Range("N53").Value = -2435.5
$ws.Range("H76").Value = 2472779.5
$ws.Range("I76").Value = 2649009.8
$ws.Range("K76").Value = 2649009.8
$ws.Range("M76").Value = -2648694.8
$ws.Range("H79").Value = 2472779.5
$ws.Range("I79").Value = 2649009.8
$ws.Range("K79").Value = 2649009.8
$ws.Range("M79").Value = -2647917.8
$ws.Range("H86").Value = 2076.7778
$ws.Range("I86").Value = 1557.4286
$ws.Range("J86").Value = 3894.5
$ws.Range("K86").Value = 1557.4286
$ws.Range("L86").Value = 3894.5
$ws.Range("M86").Value = -434.4286
$ws.Range("N86").Value = -6140.5
$ws.Range("H89").Value = 2076.7778
$ws.Range("I89").Value = 1557.4286
$ws.Range("J89").Value = 3894.5
$ws.Range("K89").Value = 7787.143
$ws.Range("L89").Value = 19472.5
$ws.Range("M89").Value = -2171.143
$ws.Range("N89").Value = -30704.5
$ws.Range("H106").Value = 1957.1875
$ws.Range("I106").Value = 1957.1875
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1957.1875
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1326.1875
$ws.Range("N106").ClearContents()
$ws.Range("H138").Value = 2554.6562
$ws.Range("I138").Value = 1751.1765
$ws.Range("J138").Value = 3465.2666
$ws.Range("K138").Value = 5253.529500000001
$ws.Range("L138").Value = 10395.7998
$ws.Range("M138").Value = -113.5295000000006
$ws.Range("N138").Value = -20675.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 319378.8
$ws.Range("I32").Value = 4906.6343
$ws.Range("J32").Value = 3184569.8
$ws.Range("K32").Value = 4906.6343
$ws.Range("L32").Value = 3184569.8
$ws.Range("M32").Value = -4619.6343
$ws.Range("N32").Value = -3185143.8
$ws.Range("H132").Value = 2524.9707
$ws.Range("I132").Value = 2193.4443
$ws.Range("J132").Value = 3803.7144
$ws.Range("K132").Value = 6580.3329
$ws.Range("L132").Value = 11411.1432
$ws.Range("M132").Value = -4050.3329
$ws.Range("N132").Value = -16471.1432
$ws.Range("H138").Value = 60516.125
$ws.Range("J138").Value = 60516.125
$ws.Range("L138").Value = 60516.125
$ws.Range("N138").Value = -70796.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2372.0889
$ws.Range("I20").Value = 1822.2307
$ws.Range("J20").Value = 3124.5264
$ws.Range("K20").Value = 1822.2307
$ws.Range("L20").Value = 3124.5264
$ws.Range("M20").Value = -1575.2307
$ws.Range("N20").Value = -3618.5264
$ws.Range("H107").Value = 4168.647
$ws.Range("I107").Value = 4243.0625
$ws.Range("J107").Value = 2978
$ws.Range("K107").Value = 4243.0625
$ws.Range("L107").Value = 2978
$ws.Range("M107").Value = -2323.0625
$ws.Range("N107").Value = -6818
$ws.Range("H134").Value = 6017.4688
$ws.Range("I134").Value = 879.96295
$ws.Range("J134").Value = 33760
$ws.Range("K134").Value = 2639.88885
$ws.Range("L134").Value = 101280
$ws.Range("M134").Value = -104.8888499999998
$ws.Range("N134").Value = -106350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 777.9583
$ws.Range("I16").Value = 799.6316
$ws.Range("J16").Value = 695.6
$ws.Range("K16").Value = 799.6316
$ws.Range("L16").Value = 695.6
$ws.Range("M16").Value = -512.6316
$ws.Range("N16").Value = -1269.6
$ws.Range("H31").Value = 4530.6553
$ws.Range("I31").Value = 3139.5833
$ws.Range("J31").Value = 5512.5884
$ws.Range("K31").Value = 3139.5833
$ws.Range("L31").Value = 5512.5884
$ws.Range("M31").Value = -2844.5833
$ws.Range("N31").Value = -6102.5884
$ws.Range("H34").Value = 4530.6553
$ws.Range("I34").Value = 3139.5833
$ws.Range("J34").Value = 5512.5884
$ws.Range("K34").Value = 3139.5833
$ws.Range("L34").Value = 5512.5884
$ws.Range("M34").Value = -2937.5833
$ws.Range("N34").Value = -5916.5884
$ws.Range("H60").Value = 7939
$ws.Range("J60").Value = 8150.5
$ws.Range("L60").Value = 8150.5
$ws.Range("N60").Value = -9172.5
$ws.Range("H68").Value = 16080.2
$ws.Range("J68").Value = 16080.2
$ws.Range("L68").Value = 16080.2
$ws.Range("N68").Value = -17578.2
$ws.Range("H71").Value = 16080.2
$ws.Range("J71").Value = 16080.2
$ws.Range("L71").Value = 48240.60000000001
$ws.Range("N71").Value = -55728.60000000001
$ws.Range("H74").Value = 17808.4
$ws.Range("J74").Value = 17808.4
$ws.Range("L74").Value = 17808.4
$ws.Range("N74").Value = -19556.4
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 17808.4
$ws.Range("J77").Value = 17808.4
$ws.Range("L77").Value = 53425.2
$ws.Range("N77").Value = -62161.2
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H113").Value = 777.9583
$ws.Range("I113").Value = 799.6316
$ws.Range("J113").Value = 695.6
$ws.Range("K113").Value = 799.6316
$ws.Range("L113").Value = 695.6
$ws.Range("M113").Value = 1370.3684
$ws.Range("N113").Value = -5035.6
$ws.Range("H134").Value = 1474.6666
$ws.Range("I134").Value = 1474.6666
$ws.Range("K134").Value = 4423.9998
$ws.Range("M134").Value = -1888.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11035741
$ws.Range("I70").Value = 19402624
$ws.Range("J70").Value = 6668.0454
$ws.Range("K70").Value = 19402624
$ws.Range("L70").Value = 6668.0454
$ws.Range("M70").Value = -19402354
$ws.Range("N70").Value = -7208.0454
$ws.Range("H73").Value = 11035741
$ws.Range("I73").Value = 19402624
$ws.Range("J73").Value = 6668.0454
$ws.Range("K73").Value = 19402624
$ws.Range("L73").Value = 6668.0454
$ws.Range("M73").Value = -19401688
$ws.Range("N73").Value = -8540.045399999999
$ws.Range("H80").Value = 64517.055
$ws.Range("I80").Value = 103527.91
$ws.Range("J80").Value = 3214.2856
$ws.Range("K80").Value = 103527.91
$ws.Range("L80").Value = 3214.2856
$ws.Range("M80").Value = -102529.91
$ws.Range("N80").Value = -5210.2856
$ws.Range("H83").Value = 64517.055
$ws.Range("I83").Value = 103527.91
$ws.Range("J83").Value = 3214.2856
$ws.Range("K83").Value = 517639.55
$ws.Range("L83").Value = 16071.428
$ws.Range("M83").Value = -512647.55
$ws.Range("N83").Value = -26055.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2368.8718
$ws.Range("I40").Value = 2106.3125
$ws.Range("J40").Value = 3569.1428
$ws.Range("K40").Value = 2106.3125
$ws.Range("L40").Value = 3569.1428
$ws.Range("M40").Value = -1970.3125
$ws.Range("N40").Value = -3841.1428
$ws.Range("H136").Value = 3664.1794
$ws.Range("I136").Value = 3252.923
$ws.Range("J136").Value = 3869.8076
$ws.Range("K136").Value = 9758.769
$ws.Range("L136").Value = 11609.4228
$ws.Range("M136").Value = -7208.769
$ws.Range("N136").Value = -16709.4228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 4999.3335
$ws.Range("J39").Value = 4999.3335
$ws.Range("L39").Value = 4999.3335
$ws.Range("N39").Value = -5825.3335
$ws.Range("H43").Value = 4500
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

